$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Rename header row (row 1) strings: _old -> _FV2310, _new -> _FV2404 ---
$ws.Range("A1").Value = "Segmentname_FV2310"
$ws.Range("B1").Value = "Segmentgruppe_FV2310"
$ws.Range("C1").Value = "Segment_FV2310"
$ws.Range("D1").Value = "Datenelement_FV2310"
$ws.Range("E1").Value = "Segment ID_FV2310"
$ws.Range("F1").Value = "Code_FV2310"
$ws.Range("G1").Value = "Qualifier_FV2310"
$ws.Range("H1").Value = "Beschreibung_FV2310"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2310"
$ws.Range("J1").Value = "Bedingung_FV2310"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2404"
$ws.Range("M1").Value = "Segmentgruppe_FV2404"
$ws.Range("N1").Value = "Segment_FV2404"
$ws.Range("O1").Value = "Datenelement_FV2404"
$ws.Range("P1").Value = "Segment ID_FV2404"
$ws.Range("Q1").Value = "Code_FV2404"
$ws.Range("R1").Value = "Qualifier_FV2404"
$ws.Range("S1").Value = "Beschreibung_FV2404"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2404"
$ws.Range("U1").Value = "Bedingung_FV2404"

# --- 2. Convert A1:U80 into an Excel Table (ListObject), preserving the
#        existing header-row formatting and avoiding an auto-generated
#        header dxf (stash formatting on a scratch row, clear, build the
#        table, then restore and discard the scratch row). ---
$headerRange = $ws.Range("A1:U1")
$scratch = $ws.Range("A100:U100")
$headerRange.Copy()
$scratch.PasteSpecial(-4122)
$headerRange.ClearFormats()

$tbl = $ws.ListObjects.Add(1, $ws.Range("A1:U80"), $null, 1)
$tbl.Name = "Table1"
$tbl.TableStyle = ""

$scratch.Copy()
$headerRange.PasteSpecial(-4122)
$ws.Rows.Item(100).Delete()

# --- 3. Freeze the top row ---
$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
